# Update cryptos list values (price and 1h volume change) per Oct 1 2024 refresh
# Uses a leading apostrophe to force text interpretation (values like "1.00" or
# "0.0000182" would otherwise be auto-converted to numbers), then resets the
# cell style back to Normal so no stray quote-prefix / text number-format style
# is left attached to the cell (matching the original plain-text cell styling).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.382.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.616.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.73%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'153.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -2.47%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.614.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -4.63%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.70%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.086.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000182"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.311.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.91%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.621.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.49%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -3.93%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'340.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.70%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.34%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'67.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +7.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000106"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.17%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'584.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +7.54%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.22%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.84%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -3.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.12%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'19.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'151.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.53%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'41.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.29%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +3.68%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'157.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'23.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0580"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.84%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.626"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.94%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.85%  "
$ws.Range("E51").Style = "Normal"
